# fix(general): fix bug date
# Rewrite the voucher rows with corrected codes / amounts / dates and
# expand the table from 6 data rows to 12 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: code, amount(E) or $null, percent(G) or $null, startedAt, expiredAt, createdAt
$rows = @(
    @("I0GOTQFEZC9FWW3Y", 30,  $null, "11/09/2022", "29/09/2022", "10/09/2022"),
    @("C63OACLO0NLWS0N0", 30,  $null, "11/09/2022", "29/09/2022", "10/09/2022"),
    @("56ZMXEMC13RMC380", 12,  $null, "07/09/2022", "30/09/2022", "07/09/2022"),
    @("AA33C0AJ336BJQA0", 12,  $null, "07/09/2022", "30/09/2022", "07/09/2022"),
    @("AAS3DWWRVOFC52TW", 12,  $null, "07/09/2022", "30/09/2022", "07/09/2022"),
    @("9HIHJDIPM50AC6NA", 12,  $null, "07/09/2022", "30/09/2022", "07/09/2022"),
    @("AF06WSE3TAJCAPQB", 56,  $null, "05/09/2022", "30/09/2022", "05/09/2022"),
    @("8IHBFU05FZAE9QEC", 56,  $null, "05/09/2022", "30/09/2022", "05/09/2022"),
    @("IA9ECZTMWPCQ03Q3", 56,  $null, "05/09/2022", "30/09/2022", "05/09/2022"),
    @("S2FC3E6BHPL3EBTX", 56,  $null, "05/09/2022", "30/09/2022", "05/09/2022"),
    @("V3BTFQE6CAKR3JOO", $null, 56,  "05/09/2022", "30/09/2022", "05/09/2022"),
    @("Y6B60ISEC5JJVK1Z", $null, 56,  "05/09/2022", "30/09/2022", "05/09/2022")
)

$r = 2
foreach ($row in $rows) {
    $code       = $row[0]
    $amount     = $row[1]
    $percent    = $row[2]
    $startedAt  = $row[3]
    $expiredAt  = $row[4]
    $createdAt  = $row[5]

    $ws.Cells.Item($r, 1).Value = $code
    $ws.Cells.Item($r, 2).Value = "ACTIVE"
    $ws.Cells.Item($r, 3).Value = "ONLINE"
    $ws.Cells.Item($r, 4).Value = "COUPON"

    if ($null -eq $amount) {
        $ws.Cells.Item($r, 5).Value = $null
        $ws.Cells.Item($r, 6).Value = $null
    } else {
        $ws.Cells.Item($r, 5).Value = $amount
        $ws.Cells.Item($r, 6).Value = "EUR"
    }

    if ($null -eq $percent) {
        $ws.Cells.Item($r, 7).Value = $null
    } else {
        $ws.Cells.Item($r, 7).Value = $percent
    }

    $ws.Cells.Item($r, 8).Value = $startedAt
    $ws.Cells.Item($r, 9).Value = $expiredAt
    $ws.Cells.Item($r, 10).Value = $createdAt

    $r = $r + 1
}
